$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) cells, forcing text type to match original inlineStr formatting
$priceUpdates = @{
    "D2" = '54.942.01'
    "D3" = '2.269.11'
    "D5" = '503.57'
    "D6" = '127.54'
    "D7" = '0.999'
    "D9" = '2.277.64'
    "D10" = '0.0973'
    "D15" = '2.672.78'
    "D16" = '54.898.65'
    "D18" = '2.260.83'
    "D20" = '4.17'
    "D21" = '310.56'
    "D22" = '6.55'
    "D23" = '0.998'
    "D24" = '59.75'
    "D27" = '7.47'
    "D28" = '171.30'
    "D29" = '6.11'
    "D34" = '17.91'
    "D35" = '0.998'
    "D39" = '36.64'
    "D41" = '0.373'
    "D42" = '134.85'
    "D44" = '4.85'
    "D45" = '256.56'
    "D51" = '16.38'
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = $origStyle
}

# Update Volume(1h) column (E) cells
$volumeUpdates = @{
    "E2" = '  +1.43%  '
    "E3" = '  -0.04%  '
    "E4" = '  +0.01%  '
    "E5" = '  +1.02%  '
    "E6" = '  -0.51%  '
    "E7" = '  +0.00%  '
    "E8" = '  +0.33%  '
    "E9" = '  +0.25%  '
    "E10" = '  +2.15%  '
    "E11" = '  +0.96%  '
    "E12" = '  +7.32%  '
    "E13" = '  +1.22%  '
    "E14" = '  +3.31%  '
    "E15" = '  +0.09%  '
    "E16" = '  +1.39%  '
    "E17" = '  +0.46%  '
    "E18" = '  -0.31%  '
    "E19" = '  +1.27%  '
    "E20" = '  +0.79%  '
    "E21" = '  +2.80%  '
    "E22" = '  +3.69%  '
    "E23" = '  -0.34%  '
    "E24" = '  -2.17%  '
    "E25" = '  -0.12%  '
    "E26" = '  +3.12%  '
    "E27" = '  +2.47%  '
    "E28" = '  +0.31%  '
    "E29" = '  +3.16%  '
    "E30" = '  +1.37%  '
    "E31" = '  +1.27%  '
    "E32" = '  +5.33%  '
    "E34" = '  +1.20%  '
    "E35" = '  -0.02%  '
    "E36" = '  +2.40%  '
    "E37" = '  -4.40%  '
    "E38" = '  +3.98%  '
    "E39" = '  +1.82%  '
    "E40" = '  +3.22%  '
    "E41" = '  +0.56%  '
    "E42" = '  +8.22%  '
    "E43" = '  +2.94%  '
    "E44" = '  +1.09%  '
    "E45" = '  +7.57%  '
    "E46" = '  +2.44%  '
    "E47" = '  +2.44%  '
    "E48" = '  -0.07%  '
    "E49" = '  +0.53%  '
    "E50" = '  +2.94%  '
    "E51" = '  +1.51%  '
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
